# Update the "Провеждане часовете в курса" schedule list on slide 3:
# shift each bullet's text down into the next paragraph (dropping the old
# date prefixes) and insert a new trailing bullet for the exam paragraph
# that previously shared a paragraph with its date prefix.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 4 ("11. 09. 2015 – изпит (тест и задача)") gains a new sibling
# paragraph after it, carrying the same bullet/paragraph formatting, so we
# can shift everything down by one without losing the last line of text.
$null = $tr.Paragraphs(4, 1).InsertAfter("`r11. 09. 2015 – изпит (тест и задача)")

# Now shift the text of each paragraph into the next one, then trim the
# leading date/prefix from the shifted copies, starting from the bottom so
# earlier paragraphs' original text is still available to copy forward.
$tr.Paragraphs(5, 1).Text = "изпит (тест и задача)"
$tr.Paragraphs(4, 1).Text = "подготовка за изпит"
$tr.Paragraphs(3, 1).Text = "занимания"
$tr.Paragraphs(2, 1).Text = "Всяка понеделник и сряда 19:30 – 22:30 "
$tr.Paragraphs(1, 1).Text = ""
